$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format to preserve literal numeric-looking strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.286.24"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.838.43"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "239.03"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "0.6243"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D8").Value = "0.07376"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "0.2885"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "24.75"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D12").Value = "1.833.92"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "4.948"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "0.00001052"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "0.6627"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "81.28"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "6.254"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "29.278.54"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "234.31"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "12.23"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "7.288"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "157.30"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "8.413"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "17.24"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "0.07120"
$ws.Range("E28").Value = "  +8.48%  "
$ws.Range("D29").Value = "1.482"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.020"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "4.023"
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").Value = "1.151"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "1.790"
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("D35").Value = "0.6987"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "2.587"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "0.01823"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.784"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.234.07"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "6.780"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "0.9471"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "1.991.60"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "101.16"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "65.24"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "6.914"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "1.683"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").Value = "8.907"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "0.3868"
$ws.Range("E51").Value = "  -2.15%  "
